$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Zeile im Spreadsheet" (column K) values - increment/decrement by 1
$ws.Range("K5").Value = 384
$ws.Range("K9").Value = 547
$ws.Range("K10").Value = 415
$ws.Range("K12").Value = 645
$ws.Range("K15").Value = 342
$ws.Range("K16").Value = 398
$ws.Range("K18").Value = 569
$ws.Range("K23").Value = 337
$ws.Range("K24").Value = 403
$ws.Range("K26").Value = 567
$ws.Range("K32").Value = 406
$ws.Range("K35").Value = 549
$ws.Range("K36").Value = 414
$ws.Range("K40").Value = 393
$ws.Range("K41").Value = 561
$ws.Range("K45").Value = 390
$ws.Range("K48").Value = 556
$ws.Range("K51").Value = 372
$ws.Range("K54").Value = 558
$ws.Range("K57").Value = 374
$ws.Range("K60").Value = 387
$ws.Range("K62").Value = 554
$ws.Range("K66").Value = 366
$ws.Range("K74").Value = 552
$ws.Range("K76").Value = 730
$ws.Range("K79").Value = 563
$ws.Range("K82").Value = 355
$ws.Range("K84").Value = 322
$ws.Range("K85").Value = 323
$ws.Range("K87").Value = 565

# Update the unit text in G15 - remove stray spaces
$ws.Range("G15").Value = "t/cap/year"
